$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order/content (row 1):
# A: CustomerCode*       (style 1)
# B: CustomerName*       (style 1)
# C: Address*            (style 1)
# D: ContactPhone        (style 1)
# E: ContactEmail        (style 1)
# F: DeliveryTerm        (style 3 -> new xf numFmtId=0 fontId=1 applyFont)
# G: PaymentTerm         (style 3)
# H: CurrencyCode*       (style 2 -> xf numFmtId=49 fontId=0 applyNumberFormat)
# I: CountryCode*        (style 0 default)

$ws.Range("A1").Value = "CustomerCode*"
$ws.Range("B1").Value = "CustomerName*"
$ws.Range("C1").Value = "Address*"
$ws.Range("D1").Value = "ContactPhone"
$ws.Range("E1").Value = "ContactEmail"
$ws.Range("F1").Value = "DeliveryTerm"
$ws.Range("G1").Value = "PaymentTerm"
$ws.Range("H1").Value = "CurrencyCode*"
$ws.Range("I1").Value = "CountryCode*"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 13.88671875
$ws.Columns.Item(2).ColumnWidth = 14.6640625
$ws.Columns.Item(3).ColumnWidth = 12.88671875
$ws.Columns.Item(4).ColumnWidth = 12.6640625
$ws.Columns.Item(5).ColumnWidth = 11.88671875
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(7).ColumnWidth = 12.6640625
$ws.Columns.Item(8).ColumnWidth = 13.6640625
$ws.Columns.Item(9).ColumnWidth = 12.44140625

# View settings
$ws.Application.ActiveWindow.Zoom = 175
$ws.Range("H7").Select()

$wb.Save()
